# runtime update (2025-10-27 20:30:04)
# Adds 2025-10-26 KHL match results to Matches_SOG, and rolls the derived
# shots-on-goal aggregates on Shots_HA / Shots_Summary / Meta_ext forward
# to reflect them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Matches_SOG — append the four games played 2025-10-26
# ---------------------------------------------------------------------
$matchesWs = $wb.Worksheets.Item("Matches_SOG")

$newMatches = @(
    @{ Row = 394; Uid = "897691"; Home = "ЦСКА";       Away = "Амур";         SogHome = 21; SogAway = 33 },
    @{ Row = 395; Uid = "897692"; Home = "Торпедо";     Away = "Автомобилист"; SogHome = 27; SogAway = 47 },
    @{ Row = 396; Uid = "897693"; Home = "Северсталь";  Away = "Драконы";      SogHome = 28; SogAway = 18 },
    @{ Row = 397; Uid = "897694"; Home = "Спартак";     Away = "ХК Сочи";      SogHome = 42; SogAway = 34 }
)

foreach ($m in $newMatches) {
    $r = $m.Row

    # Column A (uid) looks numeric ("897691") but must stay a text cell,
    # matching every other row in this column — force text via NumberFormat
    # before assigning so it isn't auto-coerced into a number.
    $uidCell = $matchesWs.Cells.Item($r, 1)
    $uidCell.NumberFormat = "@"
    $uidCell.Value = $m.Uid

    $matchesWs.Cells.Item($r, 2).Value = "2025-10-26T17:00:00"
    $matchesWs.Cells.Item($r, 3).Value = $m.Home
    $matchesWs.Cells.Item($r, 4).Value = $m.Away
    $matchesWs.Cells.Item($r, 5).Value = $m.SogHome
    $matchesWs.Cells.Item($r, 6).Value = $m.SogAway
    $matchesWs.Cells.Item($r, 7).Value = "khl_text"
}

# ---------------------------------------------------------------------
# 2) Shots_HA — as_of_utc rolls to 2025-10-26 for every team, and the
#    home/away shots-on-goal splits update for the 8 teams that played.
# ---------------------------------------------------------------------
$haWs = $wb.Worksheets.Item("Shots_HA")

for ($r = 2; $r -le 23; $r++) {
    $haWs.Cells.Item($r, 4).Value = "2025-10-26T17:00:00Z"
}

# Row => updated GP/total/per-game figures (only columns that actually move)
$haUpdates = @{
    3  = @{ F = 23; K = 643; L = 687; M = 28;   N = 29.9 }   # Автомобилист (away)
    6  = @{ F = 19; K = 549; L = 669; M = 28.9; N = 35.2 }   # Амур (away)
    10 = @{ F = 17; K = 480; L = 626; M = 28.2; N = 36.8 }   # Драконы (away)
    17 = @{ E = 13; G = 372; H = 306; I = 28.6; J = 23.5 }   # Северсталь (home)
    19 = @{ E = 23; G = 834; H = 638; I = 36.3; J = 27.7 }   # Спартак (home)
    20 = @{ E = 19; G = 603; H = 569; I = 31.7; J = 29.9 }   # Торпедо (home)
    22 = @{ F = 17; K = 434; L = 632; M = 25.5; N = 37.2 }   # ХК Сочи (away)
    23 = @{ E = 17; G = 383; H = 499; I = 22.5; J = 29.4 }   # ЦСКА (home)
}

$colIndex = @{ E = 5; F = 6; G = 7; H = 8; I = 9; J = 10; K = 11; L = 12; M = 13; N = 14 }

foreach ($rowNum in $haUpdates.Keys) {
    $rowData = $haUpdates[$rowNum]
    foreach ($col in $rowData.Keys) {
        $haWs.Cells.Item([int]$rowNum, $colIndex[$col]).Value = $rowData[$col]
    }
}

# ---------------------------------------------------------------------
# 3) Shots_Summary — same as_of_utc roll, plus combined totals for the
#    8 affected teams.
# ---------------------------------------------------------------------
$summaryWs = $wb.Worksheets.Item("Shots_Summary")

for ($r = 2; $r -le 23; $r++) {
    $summaryWs.Cells.Item($r, 4).Value = "2025-10-26T17:00:00Z"
}

$summaryUpdates = @{
    3  = @{ E = 39; F = 1115; G = 1201; H = 28.6; I = 30.8 }   # Автомобилист
    6  = @{ E = 35; F = 1039; G = 1233; H = 29.7; I = 35.2 }   # Амур
    10 = @{ E = 35; F = 984;  G = 1266; H = 28.1; I = 36.2 }   # Драконы
    17 = @{ E = 35; F = 1094; G = 886;  H = 31.3; I = 25.3 }   # Северсталь
    19 = @{ E = 35; F = 1240; G = 1086; H = 35.4; I = 31.0 }   # Спартак
    20 = @{ E = 43; F = 1451; G = 1357; H = 33.7; I = 31.6 }   # Торпедо
    22 = @{ E = 33; F = 916;  G = 1134; H = 27.8; I = 34.4 }   # ХК Сочи
    23 = @{ E = 35; F = 835;  G = 1029; I = 29.4 }             # ЦСКА (H unchanged: 23.9)
}

$summaryColIndex = @{ E = 5; F = 6; G = 7; H = 8; I = 9 }

foreach ($rowNum in $summaryUpdates.Keys) {
    $rowData = $summaryUpdates[$rowNum]
    foreach ($col in $rowData.Keys) {
        $summaryWs.Cells.Item([int]$rowNum, $summaryColIndex[$col]).Value = $rowData[$col]
    }
}

# ---------------------------------------------------------------------
# 4) Meta_ext — bump as_of_utc and the build_version counter.
# ---------------------------------------------------------------------
$metaWs = $wb.Worksheets.Item("Meta_ext")
$metaWs.Range("B2").Value = "2025-10-26T17:00:00Z"
$metaWs.Range("D2").Value = 9
